# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" timestamps and "Priority" flags for
# the rows whose Status is "Ready for handoff" across the Overview, zh-cn and
# de-de worksheets.

$wb = $excel.ActiveWorkbook

# Rows (in each table) corresponding to files that were just (re)handed off.
$rows = @(7, 8, 9, 10, 12, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-26 10:19:36"
}

# --- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-26 10:19:31"
}

# --- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-26 10:19:36"
}
